$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(9).Delete()
$ws.Range("A1").Value = "Mã"
$ws.Range("B1").Value = "Tên Sản Phẩm"
$ws.Range("C1").Value = "Hình Ảnh"
$ws.Range("D1").Value = "Giá Nhập"
$ws.Range("E1").Value = "Giá Bán"
$ws.Range("F1").Value = "Số Lượng"
$ws.Range("G1").Value = "Mô Tả"
$ws.Range("H1").Value = "Loại Hàng"
$ws.Range("K2").Select()
